# Update "horarios-141" tracking workbook with the latest scrape batch
# (31/12/2025 17:31:45-17:31:55) across the three line sheets:
#   LP1912, LP1912-215, 6203-6173
#
# For every sheet this:
#   1. Refreshes the "Última actualización: ..." timestamp in A2
#   2. Refreshes the "Total filas: N" counter in A3
#   3. Appends the new scraped rows at the bottom of the sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912  (columns: A Fecha, B Hora_Scrap, C Hora_Llegada,
#                    D Linea, E Minutos, F Parada, G Fecha)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 17:31:55"
$ws1.Cells.Item(3, 1).Value = "Total filas: 1162"

$lp1912Rows = @(
    @("17:31:45", "17:34", "10_OLMOS",            3, "LP1912", "31/12/2025"),
    @("17:31:45", "17:35", "16_P MOR-SANTA ANA",  4, "LP1912", "31/12/2025"),
    @("17:31:45", "17:38", "17X38_ROMERO",        7, "LP1912", "31/12/2025"),
    @("17:31:45", "17:47", "16_SANTA ANA",       16, "LP1912", "31/12/2025"),
    @("17:31:45", "17:50", "215_EL PELIGRO",     19, "LP1912", "31/12/2025"),
    @("17:31:45", "17:54", "10_OLMOS",           23, "LP1912", "31/12/2025"),
    @("17:31:45", "17:59", "16_SANTA ANA",       28, "LP1912", "31/12/2025"),
    @("17:31:45", "18:03", "23_HERNANDEZ",       32, "LP1912", "31/12/2025"),
    @("17:31:45", "18:04", "14_ABASTO",          33, "LP1912", "31/12/2025"),
    @("17:31:45", "18:21", "16_SANTA ANA",       50, "LP1912", "31/12/2025"),
    @("17:31:45", "18:24", "11_ETCHEVERRY",      53, "LP1912", "31/12/2025"),
    @("17:31:45", "18:30", "23_HERNANDEZ",       59, "LP1912", "31/12/2025"),
    @("17:31:45", "18:34", "14X44_ABASTO",       63, "LP1912", "31/12/2025"),
    @("17:31:45", "18:41", "16_P MOR-SANTA ANA", 70, "LP1912", "31/12/2025"),
    @("17:31:45", "18:41", "14_ABASTO",          70, "LP1912", "31/12/2025"),
    @("17:31:45", "18:51", "15_ABASTO",          80, "LP1912", "31/12/2025"),
    @("17:31:45", "18:59", "23_HERNANDEZ",       88, "LP1912", "31/12/2025"),
    @("17:31:45", "19:01", "17_ROMERO",          90, "LP1912", "31/12/2025"),
    @("17:31:45", "19:06", "14_ABASTO",          95, "LP1912", "31/12/2025")
)

$r = 1145
foreach ($row in $lp1912Rows) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215  (columns: A Fecha(blank), B Fecha, C Hora_Scrap,
#                        D Hora_Llegada, E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 17:31:55"
$ws2.Cells.Item(3, 1).Value = "Total filas: 77"

$ws2.Cells.Item(78, 2).Value = "31/12/2025"
$ws2.Cells.Item(78, 3).Value = "17:31:45"
$ws2.Cells.Item(78, 4).Value = "17:50"
$ws2.Cells.Item(78, 5).Value = "215_EL PELIGRO"
$ws2.Cells.Item(78, 6).Value = 19
$ws2.Cells.Item(78, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173  (columns: A Fecha(blank), B Fecha, C Hora_Scrap,
#                       D Hora_Llegada, E Linea, F Minutos, G Parada)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 31/12/2025 17:31:55"
$ws3.Cells.Item(3, 1).Value = "Total filas: 139"

$ws3.Cells.Item(140, 2).Value = "31/12/2025"
$ws3.Cells.Item(140, 3).Value = "17:31:50"
$ws3.Cells.Item(140, 4).Value = "18:21"
$ws3.Cells.Item(140, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(140, 6).Value = 50
$ws3.Cells.Item(140, 7).Value = "L6203"
